# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-18 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-19 Sunday", 2)

# Update the division-problem answers in the table, cell by cell so that
# identical old/new values across different cells cannot clash with each
# other during replacement.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="39÷3=13, 0"},
    @{Row=1;  Col=2; New="18÷8=2, 2"},
    @{Row=1;  Col=3; New="41÷4=10, 1"},
    @{Row=1;  Col=4; New="29÷5=5, 4"},
    @{Row=1;  Col=5; New="38÷8=4, 6"},

    @{Row=5;  Col=1; New="38÷4=9, 2"},
    @{Row=5;  Col=2; New="13÷3=4, 1"},
    @{Row=5;  Col=3; New="79÷5=15, 4"},
    @{Row=5;  Col=4; New="34÷8=4, 2"},
    @{Row=5;  Col=5; New="31÷3=10, 1"},

    @{Row=9;  Col=1; New="30÷9=3, 3"},
    @{Row=9;  Col=2; New="79÷3=26, 1"},
    @{Row=9;  Col=3; New="27÷9=3, 0"},
    @{Row=9;  Col=4; New="62÷3=20, 2"},
    @{Row=9;  Col=5; New="18÷9=2, 0"},

    @{Row=13; Col=1; New="64÷9=7, 1"},
    @{Row=13; Col=2; New="22÷2=11, 0"},
    @{Row=13; Col=3; New="64÷6=10, 4"},
    @{Row=13; Col=4; New="57÷9=6, 3"},
    @{Row=13; Col=5; New="18÷8=2, 2"},

    @{Row=17; Col=1; New="98÷7=14, 0"},
    @{Row=17; Col=2; New="73÷6=12, 1"},
    @{Row=17; Col=3; New="37÷4=9, 1"},
    @{Row=17; Col=4; New="84÷8=10, 4"},
    @{Row=17; Col=5; New="44÷3=14, 2"}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $r = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters from the range
    # so only the visible text is replaced.
    $r.End = $r.End - 1
    $r.Text = $u.New
}
